# Convert the two Word field codes (m:'...'.setDocumentSubject() and
# m:''.getDocumentSubject()) into literal M2Doc-style text runs
# "{m:...}" made up of one <w:r><w:t>...</w:t></w:r> run per original
# instrText chunk (matching TokenIteratorFieldRewriterSplit output).

$d = $word.ActiveDocument

function Insert-PlainRuns($paragraphIndex, $chunks) {
    # $chunks is an array of 2-element arrays: @(text, preserveSpace)
    $para = $d.Paragraphs($paragraphIndex)
    $pos = $para.Range.Start
    $insertionPoint = $d.Range($pos, $pos)

    $runsXml = ""
    foreach ($chunk in $chunks) {
        $text = $chunk[0]
        $preserve = $chunk[1]

        $escaped = $text -replace '&', '&amp;'
        $escaped = $escaped -replace '<', '&lt;'
        $escaped = $escaped -replace '>', '&gt;'

        if ($preserve) {
            $runsXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
        } else {
            $runsXml += '<w:r><w:t>' + $escaped + '</w:t></w:r>'
        }
    }

    $package = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body><w:p>' + $runsXml + '</w:p></w:body>' `
        + '</w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'

    $insertionPoint.InsertXML($package)
}

# --- Field 1: m:'Some value'.setDocumentSubject() -------------------------
$field1 = $d.Fields(1)
$field1ParaIndex = $field1.Code.Paragraphs(1).Index
$field1.Delete()

Insert-PlainRuns $field1ParaIndex @(
    , @("{m:", $false)
    , @("'", $false)
    , @("Some value", $false)
    , @("'", $false)
    , @(".", $false)
    , @("setDocument", $false)
    , @("Subject", $false)
    , @("()}", $true)
)

# --- Field 2: m:''.getDocumentSubject() ------------------------------------
$field2 = $d.Fields(1)
$field2ParaIndex = $field2.Code.Paragraphs(1).Index
$field2.Delete()

Insert-PlainRuns $field2ParaIndex @(
    , @("{m:''.g", $false)
    , @("etDocument", $false)
    , @("Subject", $false)
    , @("()}", $true)
)

Write-Host "Done. Doc text:" $d.Content.Text
